$d = $word.ActiveDocument

# Locate the paragraph that holds the XSS test string (starts with "Test" and
# also contains the "svg" payload). We scan by content instead of a hard-coded
# index so the script keeps working if paragraphs ever get reordered.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Test") -and $t.Contains("svg") -and $t.Contains("onload")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the target paragraph"
}

$rng = $target.Range

# Pull the paragraph's own <w:p ...> opening tag (with its w14:paraId / rsid
# attributes etc.) straight out of the live document so we reuse it verbatim
# instead of inventing new identifiers.
$openXml = $rng.WordOpenXML
if ($openXml -notmatch '(<w:p [^>]*>)') {
    throw "Could not find the paragraph's opening tag"
}
$pOpenTag = $matches[1]

# Rebuild the paragraph contents:
#  - "Test" + curly-quote pair + "><" becomes three separate runs, with the
#    curly quotes normalized to straight ' and " characters.
#  - the svg/onload payload keeps its spellStart/spellEnd proofErr markers
#    around "svg", but the three runs that used to be split apart by the now
#    removed gramStart/gramEnd proofErr markers are merged back into one run.
$newParagraphInner = (
    '<w:r><w:t>Test</w:t></w:r>' +
    '<w:r><w:t>' + [char]0x27 + '&quot;</w:t></w:r>' +
    '<w:r><w:t>&gt;&lt;</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>svg</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>/onload=alert(1)/&gt;</w:t></w:r>'
)

$packageXml = (
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body>' +
    $pOpenTag + $newParagraphInner + '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
)

[void]$rng.InsertXML($packageXml)
